$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule row (E8): "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active selection moving to E8, as recorded in the saved file
$ws.Activate()
$ws.Range("E8").Select()
